$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 177, shifting the existing rows 177-192 down to 178-193.
$ws.Rows.Item(177).Insert()

# Populate the newly inserted row 177 with the new weekly price record.
$ws.Range("A177").Value = 7
$ws.Range("B177").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C177").Value = "Ñuble"
$ws.Range("D177").Value = 44610
$ws.Range("E177").Value = 16
$ws.Range("F177").Value = 100112032
$ws.Range("G177").Value = "Zapallo italiano"
$ws.Range("H177").Value = "Sin especificar"
$ws.Range("I177").Value = "Primera"
$ws.Range("J177").Value = 100
$ws.Range("K177").Value = 6000
$ws.Range("L177").Value = 6500
$ws.Range("M177").Value = 6250
$ws.Range("N177").Value = "$/caja 50 unidades"
$ws.Range("O177").Value = "Región del Maule"
$ws.Range("P177").Value = 125
$ws.Range("Q177").Value = 50
$ws.Range("R177").Value = "Hortaliza"
